$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.532
$ws.Range("A8").Value = -22.36460000000002
$ws.Range("A10").Value = -21.73489999999999
$ws.Range("A12").Value = -21.8769
$ws.Range("C12").Value = -13.22989999999999
$ws.Range("D12").Value = -8.566299999999995
$ws.Range("D13").Value = -8.988199999999985
$ws.Range("C15").Value = -14.24309999999999
$ws.Range("C17").Value = -13.5636
$ws.Range("A18").Value = -22.1087
$ws.Range("D21").Value = -8.120999999999995
$ws.Range("D25").Value = -7.350900000000001
$ws.Range("C26").Value = -12.5928
$ws.Range("C27").Value = -13.21099999999999
$ws.Range("C28").Value = -13.6713
$ws.Range("D32").Value = -6.639999999999999
$ws.Range("D36").Value = -7.740999999999999
$ws.Range("A37").Value = -20.02719999999999
$ws.Range("C37").Value = -12.9294
$ws.Range("D38").Value = -7.883499999999998
$ws.Range("D41").Value = -8.224999999999998
$ws.Range("C47").Value = -12.7604
$ws.Range("D52").Value = -7.900900000000004
$ws.Range("A55").Value = -22.3171
$ws.Range("D59").Value = -8.333999999999993
$ws.Range("C65").Value = -12.2852
$ws.Range("D67").Value = -7.023899999999997
$ws.Range("A68").Value = -21.47249999999999
$ws.Range("C73").Value = -11.08090000000001
$ws.Range("A77").Value = -20.43709999999999
$ws.Range("A78").Value = -20.09209999999997
$ws.Range("A81").Value = -22.01840000000001
$ws.Range("A82").Value = -21.8587
$ws.Range("C84").Value = -13.72449999999999
$ws.Range("D84").Value = -7.964399999999991
$ws.Range("C85").Value = -12.9826
$ws.Range("D88").Value = -7.929799999999993
$ws.Range("D89").Value = -8.199699999999993
$ws.Range("C93").Value = -10.2461
$ws.Range("C95").Value = -13.31589999999999
$ws.Range("D95").Value = -7.633999999999999
$ws.Range("C98").Value = -13.19830000000001
$ws.Range("C99").Value = -12.21980000000001
$ws.Range("C101").Value = -13.00930000000001
$ws.Range("D105").Value = -8.382300000000004
